$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.719.69"
$ws.Range("E2").Value = "  +4.59%  "
$ws.Range("D3").Value = "'2.659.19"
$ws.Range("E3").Value = "  +3.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'567.55"
$ws.Range("E5").Value = "  +6.25%  "
$ws.Range("D6").Value = "'146.41"
$ws.Range("E6").Value = "  +3.74%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.608"
$ws.Range("E8").Value = "  +3.94%  "
$ws.Range("D9").Value = "'2.658.13"
$ws.Range("E9").Value = "  +3.15%  "
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("E11").Value = "  +5.65%  "
$ws.Range("E12").Value = "  +6.89%  "
$ws.Range("E13").Value = "  +4.39%  "
$ws.Range("D14").Value = "'3.113.69"
$ws.Range("E14").Value = "  +2.66%  "
$ws.Range("D15").Value = "'60.638.94"
$ws.Range("E15").Value = "  +4.56%  "
$ws.Range("D16").Value = "'22.08"
$ws.Range("E16").Value = "  +6.75%  "
$ws.Range("E17").Value = "  +5.51%  "
$ws.Range("D18").Value = "'2.656.59"
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("D19").Value = "'4.54"
$ws.Range("E19").Value = "  +3.52%  "
$ws.Range("D20").Value = "'342.70"
$ws.Range("E21").Value = "  +4.36%  "
$ws.Range("E22").Value = "  +3.73%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'66.36"
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("D25").Value = "'0.438"
$ws.Range("E25").Value = "  +4.54%  "
$ws.Range("D26").Value = "'0.165"
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'7.40"
$ws.Range("E28").Value = "  +5.45%  "
$ws.Range("D29").Value = "'0.0₃0805"
$ws.Range("E29").Value = "  +11.26%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("E31").Value = "  +5.09%  "
$ws.Range("D32").Value = "'6.16"
$ws.Range("E32").Value = "  +4.98%  "
$ws.Range("D33").Value = "'159.37"
$ws.Range("E33").Value = "  +2.97%  "
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("E35").Value = "  +5.94%  "
$ws.Range("D36").Value = "'0.896"
$ws.Range("E36").Value = "  +7.92%  "
$ws.Range("E37").Value = "  +5.93%  "
$ws.Range("D38").Value = "'0.890"
$ws.Range("E38").Value = "  +8.87%  "
$ws.Range("E39").Value = "  +8.02%  "
$ws.Range("E40").Value = "  +1.67%  "
$ws.Range("D41").Value = "'298.75"
$ws.Range("E41").Value = "  +6.43%  "
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +4.01%  "
$ws.Range("D45").Value = "'0.602"
$ws.Range("E45").Value = "  +2.51%  "
$ws.Range("E46").Value = "  +2.52%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'128.32"
$ws.Range("E47").Value = "  +16.59%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'19.39"
$ws.Range("E48").Value = "  +3.14%  "
$ws.Range("D49").Value = "'10.72"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("E50").Value = "  +4.29%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'18.75"
$ws.Range("E51").Value = "  +5.79%  "
